# Inventaire : ajout de 13 clés USB sur la feuille active
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Clé USB 1",
    "Clé USB 2",
    "Clé USB 3",
    "Clé USB 4",
    "Clé USB 5",
    "Clé USB 6",
    "Clé USB 7",
    "Clé USB 8",
    "Clé USB 9",
    "Clé USB 10",
    "Clé USB 11",
    "Clé USB 12",
    "Clé USB 13"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# La feuille se retrouve avec la sélection sur G14 après la saisie
$ws.Range("G14").Select()
